$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (Sheet1) and name it Sheet2
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Populate the new sheet with header + data rows
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "John Doe"
$ws2.Range("B2").Value = "RyzK/uu8Q18DVx3DwtARbQ=="
$ws2.Range("A3").Value = "John Doe"
$ws2.Range("B3").Value = "g3/DOGG74jC3Flrr3yH+3D/yKbOqqUNM"

# Autofit column B like the original sheet's bestFit columns
$ws2.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Select B3 on the new sheet and make it the active sheet/tab
$ws2.Range("B3").Select()
$ws2.Activate()
